$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$n = 24
$identity = New-Object 'object[,]' $n, $n
for ($r = 0; $r -lt $n; $r++) {
    for ($c = 0; $c -lt $n; $c++) {
        if ($r -eq $c) {
            $identity[$r, $c] = 1.0
        } else {
            $identity[$r, $c] = 0.0
        }
    }
}

$ws.Range("A1:X24").Value = $identity
